# Update the "join-no-learning" sheet:
#  - delete the two "4o-vision" rows (row 7 under Zero-Shot, row 12 under Few-Shot)
#  - fill in corrected/ new values for the "4o-mini" rows (With Heuristics columns,
#    plus corrected FPR/FNR under No Heuristics)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("join-no-learning")

# --- Fill in the "4o-mini" row under Zero-Shot (row 5) ---
# No Heuristics: ACC unchanged, FPR/FNR corrected
$ws.Range("E5").Value = 0.17204301075268799
$ws.Range("F5").Value = 0.38318670576734998
# With Heuristics: ACC/FPR/FNR newly populated
$ws.Range("G5").Value = 0.81264255457803802
$ws.Range("H5").Value = 0.17139133268165499
$ws.Range("I5").Value = 0.015966112740306199

# --- Fill in the "4o-mini" row under Few-Shot (row 10) ---
# No Heuristics: ACC unchanged, FPR/FNR corrected
$ws.Range("E10").Value = 0.207559465623981
$ws.Range("F10").Value = 0.27891821440208497
# With Heuristics: ACC/FPR/FNR newly populated
$ws.Range("G10").Value = 0.87976539589442804
$ws.Range("H10").Value = 0.119908765070055
$ws.Range("I10").Value = 0.000325839035516454

# --- Delete the "4o-vision" row under Few-Shot (old row 12) first (higher row number) ---
$ws.Rows.Item(12).Delete()

# --- Delete the "4o-vision" row under Zero-Shot (old row 7) ---
$ws.Rows.Item(7).Delete()

# Match the final cell selection recorded in the saved workbook
$ws.Range("J15").Select() | Out-Null

$wb.Save()
